# Updated symbol list on Tue Dec 20 12:45:43 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep a purely textual value (matching the original
    # inlineStr/shared-string cells) instead of letting Excel reinterpret a
    # numeric-looking string as a real number.
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

Set-TextValue "D2" "249.13"
Set-TextValue "D3" "21.96"
Set-TextValue "D4" "5.424"
Set-TextValue "D5" "0.05639"
Set-TextValue "D8" "0.8168"
Set-TextValue "D9" "0.9182"

$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.0005761"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1439"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07474"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03200"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03087"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09326"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.557"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001596"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04731"
$ws.Range("E18").Value = "17CoinExTokenCET"

Set-TextValue "D19" "0.006375"
Set-TextValue "D20" "0.005071"
Set-TextValue "D21" "0.001031"
Set-TextValue "D23" "3.730"
Set-TextValue "D24" "2.150"
Set-TextValue "D25" "0.3308"
Set-TextValue "D26" "0.1296"
Set-TextValue "D40" "0.04013"
Set-TextValue "D41" "0.006891"
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.003400"
Set-TextValue "D44" "0.007717"
Set-TextValue "D45" "0.00005573"
Set-TextValue "D48" "0.6753"
Set-TextValue "D49" "0.2122"
